$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3421
$ws1.Range("F4").Value = 582
$ws1.Range("F5").Value = 838
$ws1.Range("F6").Value = 317
$ws1.Range("F10").Value = 639
$ws1.Range("F12").Value = 440
$ws1.Range("F13").Value = 69
$ws1.Range("G14").Value = 69.90000000000001
$ws1.Range("F15").Value = 343
$ws1.Range("F16").Value = 60
$ws1.Range("F18").Value = 98
$ws1.Range("F19").Value = 186

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F4").Value = 745
$ws3.Range("F5").Value = 1787

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 745
$ws4.Range("F5").Value = 1787
$ws4.Range("F6").Value = 3421
$ws4.Range("F9").Value = 582
$ws4.Range("F10").Value = 838
$ws4.Range("F11").Value = 317
$ws4.Range("F20").Value = 639
$ws4.Range("F24").Value = 440
$ws4.Range("F26").Value = 69
$ws4.Range("G27").Value = 69.90000000000001
$ws4.Range("F29").Value = 343
$ws4.Range("F30").Value = 60
$ws4.Range("F34").Value = 98
$ws4.Range("F40").Value = 186
